# Apply the "update splashloader und präsentation" edits to slide 2
# (Implementierungs-Statistik slide): update commit count, package count,
# and source-line count / resize its textbox to fit the new (longer) text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- "300 Commits" -> "312 Commits" -------------------------------------
# The shape holds two runs: "300 " and "Commits". Only change the digits
# run, leaving the "Commits" run (and its rPr) untouched.
$shpCommits = $s.Shapes.Item(5)
$trCommits = $shpCommits.TextFrame.TextRange
$trCommits.Characters(1, 4).Text = "312 "

# --- "20 Pakete" -> "21 Pakete" ------------------------------------------
$shpPakete = $s.Shapes.Item(7)
$shpPakete.TextFrame.TextRange.Text = "21 Pakete"

# --- "15000 Zeilen Quellcode" -> "20.924 Zeilen Quellcode" ---------------
# Originally two runs: "15000 " (plain) and "Zeilen Quellcode" (dirty="0").
# The target keeps a single run with the "dirty=0" formatting, so first
# rewrite the second run to contain the full new text (it keeps its own
# rPr), then delete the leftover text that used to be the first run.
$shpCode = $s.Shapes.Item(8)
$trCode = $shpCode.TextFrame.TextRange
$oldLen = $trCode.Length
$trCode.Characters(7, $oldLen - 6).Text = "20.924 Zeilen Quellcode"
$trCode.Characters(1, 6).Text = ""

# Resize/reposition the textbox (it grew because the new text is longer).
$shpCode.Left = 576.3699212598425
$shpCode.Width = 291.56481
